# Updated cryptos list on Tue Mar 26 02:27:20 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.250.05'
$ws.Range('E2').Value = '  +5.48%  '

# Row 3
$ws.Range('D3').Value = '3.613.32'

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = '''592.13'
$ws.Range('E5').Value = '  +3.73%  '

# Row 6
$ws.Range('D6').Value = '''190.93'
$ws.Range('E6').Value = '  +4.03%  '

# Row 7
$ws.Range('D7').Value = '''0.644'
$ws.Range('E7').Value = '  +2.15%  '

# Row 8
$ws.Range('D8').Value = '3.607.44'
$ws.Range('E8').Value = '  +5.29%  '

# Row 9
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  -0.06%  '

# Row 10
$ws.Range('E10').Value = '  +3.79%  '

# Row 11
$ws.Range('E11').Value = '  +3.65%  '

# Row 12
$ws.Range('D12').Value = '''58.66'
$ws.Range('E12').Value = '  +4.62%  '

# Row 13
$ws.Range('D13').Value = '''0.0000290'
$ws.Range('E13').Value = '  +4.76%  '

# Row 14
$ws.Range('D14').Value = '''9.89'
$ws.Range('E14').Value = '  +5.62%  '

# Row 15
$ws.Range('D15').Value = '4.186.04'
$ws.Range('E15').Value = '  +5.19%  '

# Row 16
$ws.Range('D16').Value = '''19.68'
$ws.Range('E16').Value = '  +6.18%  '

# Row 17
$ws.Range('D17').Value = '3.610.60'
$ws.Range('E17').Value = '  +5.09%  '

# Row 18
$ws.Range('D18').Value = '70.217.55'
$ws.Range('E18').Value = '  +5.52%  '

# Row 19
$ws.Range('D19').Value = '''12.58'
$ws.Range('E19').Value = '  +4.55%  '

# Row 20
$ws.Range('E20').Value = '  +0.70%  '

# Row 21
$ws.Range('D21').Value = '''1.06'
$ws.Range('E21').Value = '  +4.35%  '

# Row 22
$ws.Range('D22').Value = '''489.78'
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').Value = '''19.45'
$ws.Range('E23').Value = '  +17.02%  '

# Row 24
$ws.Range('D24').Value = '''5.36'
$ws.Range('E24').Value = '  +6.74%  '

# Row 25
$ws.Range('D25').Value = '''4.48'
$ws.Range('E25').Value = '  +6.09%  '

# Row 26
$ws.Range('D26').Value = '''90.97'
$ws.Range('E26').Value = '  +2.07%  '

# Row 27
$ws.Range('D27').Value = '''3.13'
$ws.Range('E27').Value = '  +6.24%  '

# Row 28
$ws.Range('D28').Value = '''11.22'
$ws.Range('E28').Value = '  +1.55%  '

# Row 29
$ws.Range('D29').Value = '''9.56'
$ws.Range('E29').Value = '  +4.49%  '

# Row 30
$ws.Range('D30').Value = '''32.90'
$ws.Range('E30').Value = '  +5.07%  '

# Row 31
$ws.Range('D31').Value = '''7.68'
$ws.Range('E31').Value = '  +7.11%  '

# Row 32
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '''629.82'
$ws.Range('E32').Value = '  +6.29%  '

# Row 33
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''12.33'
$ws.Range('E33').Value = '  +5.74%  '

# Row 34
$ws.Range('E34').Value = '  +7.03%  '

# Row 35
$ws.Range('D35').Value = '''65.89'
$ws.Range('E35').Value = '  +4.35%  '

# Row 36
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '''38.83'
$ws.Range('E36').Value = '  +7.66%  '

# Row 37
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0819'
$ws.Range('E37').Value = '  +6.65%  '

# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').Value = '''0.404'
$ws.Range('E38').Value = '  +5.48%  '

# Row 39
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.03%  '

# Row 40
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.146'
$ws.Range('E40').Value = '  -0.67%  '

# Row 41
$ws.Range('D41').Value = '''3.56'
$ws.Range('E41').Value = '  -0.65%  '

# Row 42
$ws.Range('D42').Value = '3.312.72'
$ws.Range('E42').Value = '  +4.84%  '

# Row 43
$ws.Range('E43').Value = '  +7.63%  '

# Row 44
$ws.Range('D44').Value = '''0.0453'
$ws.Range('E44').Value = '  +6.02%  '

# Row 45
$ws.Range('D45').Value = '''2.72'
$ws.Range('E45').Value = '  +7.42%  '

# Row 46
$ws.Range('E46').Value = '  +2.98%  '

# Row 47
$ws.Range('E47').Value = '  +3.06%  '

# Row 48
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '''2.73'
$ws.Range('E48').Value = '  -2.33%  '

# Row 49
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '''9.06'
$ws.Range('E49').Value = '  +3.66%  '

# Row 50
$ws.Range('D50').Value = '''3.29'
$ws.Range('E50').Value = '  +3.84%  '

# Strip the quote-prefix style COM adds for the text-forced price cells above,
# so number-look-alike prices stay plain General-formatted text cells.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'

